# Updates cryptos list values (Price column D, Volume(1h) column E)
# per the commit "Updated cryptos list on Tue Oct 22 07:40:40 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.532.95'
$ws.Range('E2').Value = '  -1.62%  '
$ws.Range('D3').Value = '2.653.32'
$ws.Range('E3').Value = '  -2.75%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.69'
$ws.Range('E5').Value = '  -1.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.76'
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.548'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = '2.653.18'
$ws.Range('E9').Value = '  -2.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.145'
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.366'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.12'
$ws.Range('E14').Value = '  -1.97%  '
$ws.Range('D15').Value = '3.137.32'
$ws.Range('E15').Value = '  -2.72%  '
$ws.Range('E16').Value = '  -3.16%  '
$ws.Range('D17').Value = '67.574.71'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '2.653.08'
$ws.Range('E18').Value = '  -2.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.16'
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.20'
$ws.Range('E20').Value = '  +7.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '363.56'
$ws.Range('E21').Value = '  -2.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.41'
$ws.Range('E22').Value = '  -2.25%  '
$ws.Range('E23').Value = '  -3.64%  '
$ws.Range('E25').Value = '  -4.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.28'
$ws.Range('E26').Value = '  -3.16%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = '2.791.22'
$ws.Range('E28').Value = '  -2.64%  '
$ws.Range('E29').Value = '  -2.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '558.58'
$ws.Range('E31').Value = '  -5.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.06'
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('E33').Value = '  -3.30%  '
$ws.Range('E34').Value = '  -3.19%  '
$ws.Range('E35').Value = '  +2.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.55'
$ws.Range('E37').Value = '  -4.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '157.63'
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.40'
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('E40').Value = '  -1.92%  '
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('E42').Value = '  -3.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.92'
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('E44').Value = '  -4.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.30'
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('D47').Value = '0.0₆0300'
$ws.Range('E47').Value = '  -3.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.596'
$ws.Range('E48').Value = '  -1.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '154.59'
$ws.Range('E49').Value = '  -0.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.88'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('E51').Value = '  -2.86%  '
